$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "C2" = 84.32042520880789
    "D2" = 78.40909090909091
    "E2" = 86.79245283018868
    "F2" = 82.38805970149254

    "C3" = 84.66211085801064
    "D3" = 83.33333333333334
    "E3" = 79.58633093525181
    "F3" = 81.41674333026678

    "C4" = 84.43432042520881
    "D4" = 84.34442270058709
    "E4" = 77.51798561151078
    "F4" = 80.78725398313027

    "C5" = 83.44722854973425
    "D5" = 80.67150635208712
    "E5" = 79.94604316546763
    "F5" = 80.30713640469737

    "C6" = 82.75731105203189
    "D6" = 78.31325301204819
    "E6" = 81.83453237410072
    "F6" = 80.03518029903253
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
